$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3, shifting existing rows 3..85 down to 4..86
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with the new record
$ws.Cells.Item(3, 1).Value = 7
$ws.Cells.Item(3, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(3, 3).Value = "Ñuble"
$ws.Cells.Item(3, 4).Value = 44691
$ws.Cells.Item(3, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(3, 5).Value = 16
$ws.Cells.Item(3, 6).Value = 100112021
$ws.Cells.Item(3, 7).Value = "Ají"
$ws.Cells.Item(3, 8).Value = "Cristal"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 60
$ws.Cells.Item(3, 11).Value = 24000
$ws.Cells.Item(3, 12).Value = 25000
$ws.Cells.Item(3, 13).Value = 24500
$ws.Cells.Item(3, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(3, 15).Value = "Región del Maule"
$ws.Cells.Item(3, 16).Value = 980
$ws.Cells.Item(3, 17).Value = 25
$ws.Cells.Item(3, 18).Value = "Hortaliza"
